# Auto-generated: apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "general"
}

$ws.Range("D2").Value = "67.556.02"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "2.613.81"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue "D5" "600.99"
$ws.Range("E5").Value = "  +1.93%  "
Set-TextValue "D6" "153.86"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue "D8" "0.549"
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("D9").Value = "2.611.43"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("E10").Value = "  +11.86%  "
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("E13").Value = "  +0.11%  "
Set-TextValue "D14" "27.90"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("E15").Value = "  +4.55%  "
$ws.Range("D16").Value = "3.089.65"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("D17").Value = "67.577.89"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").Value = "2.622.44"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("E19").Value = "  +0.31%  "
Set-TextValue "D20" "362.56"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  +4.72%  "
$ws.Range("E24").Value = "  +0.00%  "
Set-TextValue "D25" "69.94"
$ws.Range("E25").Value = "  +3.54%  "
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("E27").Value = "  +5.04%  "
$ws.Range("D28").Value = "2.747.40"
Set-TextValue "D29" "581.64"
$ws.Range("E29").Value = "  -1.63%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  -0.13%  "
Set-TextValue "D32" "7.92"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("E37").Value = "  +0.06%  "
Set-TextValue "D38" "19.43"
$ws.Range("E38").Value = "  +1.51%  "
Set-TextValue "D39" "155.57"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  +3.84%  "
Set-TextValue "D43" "2.66"
$ws.Range("E43").Value = "  +3.73%  "
Set-TextValue "D44" "41.09"
$ws.Range("E44").Value = "  -0.37%  "
Set-TextValue "D45" "0.999"
$ws.Range("E45").Value = "  +0.03%  "
Set-TextValue "D46" "16.42"
$ws.Range("E46").Value = "  +0.07%  "
Set-TextValue "D47" "156.70"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("E48").Value = "  -5.35%  "
$ws.Range("E49").Value = "  +0.72%  "
Set-TextValue "D50" "20.98"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").Value = "  +1.51%  "
